$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1, Col 1: 180×7= -> 612×3=
$t.Cell(1, 1).Range.Text = "612×3="
# Row 1, Col 2: 453×2= -> 528×7=
$t.Cell(1, 2).Range.Text = "528×7="
# Row 1, Col 3: 508×9= -> 653×2=
$t.Cell(1, 3).Range.Text = "653×2="
# Row 1, Col 4: 733×7= -> 469×3=
$t.Cell(1, 4).Range.Text = "469×3="
# Row 1, Col 5: 552×2= -> 143×5=
$t.Cell(1, 5).Range.Text = "143×5="

# Row 5, Col 1: 846×8= -> 503×4=
$t.Cell(5, 1).Range.Text = "503×4="
# Row 5, Col 2: 289×3= -> 562×9=
$t.Cell(5, 2).Range.Text = "562×9="
# Row 5, Col 3: 439×2= -> 376×9=
$t.Cell(5, 3).Range.Text = "376×9="
# Row 5, Col 4: 965×6= -> 774×4=
$t.Cell(5, 4).Range.Text = "774×4="
# Row 5, Col 5: 597×5= -> 982×5=
$t.Cell(5, 5).Range.Text = "982×5="

# Row 10, Col 1: 125×4= -> 483×5=
$t.Cell(10, 1).Range.Text = "483×5="
# Row 10, Col 2: 780×3= -> 110×3=
$t.Cell(10, 2).Range.Text = "110×3="
# Row 10, Col 3: 861×4= -> 966×6=
$t.Cell(10, 3).Range.Text = "966×6="
# Row 10, Col 4: 884×9= -> 113×9=
$t.Cell(10, 4).Range.Text = "113×9="
# Row 10, Col 5: 558×8= -> 660×8=
$t.Cell(10, 5).Range.Text = "660×8="

# Row 15, Col 1: 975×7= -> 217×6=
$t.Cell(15, 1).Range.Text = "217×6="
# Row 15, Col 2: 622×3= -> 697×2=
$t.Cell(15, 2).Range.Text = "697×2="
# Row 15, Col 3: 707×9= -> 848×7=
$t.Cell(15, 3).Range.Text = "848×7="
# Row 15, Col 4: 921×7= -> 538×2=
$t.Cell(15, 4).Range.Text = "538×2="
# Row 15, Col 5: 180×7= -> 848×3=
$t.Cell(15, 5).Range.Text = "848×3="

# Row 20, Col 1: 468×8= -> 513×4=
$t.Cell(20, 1).Range.Text = "513×4="
# Row 20, Col 2: 567×2= -> 700×7=
$t.Cell(20, 2).Range.Text = "700×7="
# Row 20, Col 3: 906×3= -> 945×4=
$t.Cell(20, 3).Range.Text = "945×4="
# Row 20, Col 4: 832×8= -> 499×8=
$t.Cell(20, 4).Range.Text = "499×8="
# Row 20, Col 5: 588×8= -> 905×9=
$t.Cell(20, 5).Range.Text = "905×9="
